# Daily attendance processing - 2025-10-17 10:21:09
# For every "Recorded By" (column G) cell that lists multiple recorders
# separated by ", ", reverse the order of the listed names/emails.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2

    if ($val -ne $null -and $val -is [string] -and $val.Contains(",")) {
        $parts = $val -split ", "
        $n = $parts.Count

        $rev = @()
        for ($i = $n - 1; $i -ge 0; $i--) {
            $rev += $parts[$i]
        }

        $joined = [string]::Join(", ", $rev)
        $cell.Value = $joined
    }
}
